# Regenerate the handback-status report for the new handback batch:
#   2072ad18-79c4-4d69-b576-ad39c016a329  -> e1ef11d8-2e75-4ed4-b1d4-b8e55d6bb396
#   e59dc089-a281-4a81-8204-8859527d7dee  -> ffffd7ece1c6-aa6f-4143-b112-7788f6378cb3
# plus refreshed handoff/handback timestamps and xliff correspondence file names.

$wb = $excel.ActiveWorkbook

$oldFile1 = "2072ad18-79c4-4d69-b576-ad39c016a329.md"
$oldFile2 = "e59dc089-a281-4a81-8204-8859527d7dee.md"
$newFile1 = "e1ef11d8-2e75-4ed4-b1d4-b8e55d6bb396.md"
$newFile2 = "ffffd7ece1c6-aa6f-4143-b112-7788f6378cb3.md"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFile1
$wsOverview.Range("G2").Value = "2016-08-28 05:05:41"
$wsOverview.Range("A3").Value = $newFile2
$wsOverview.Range("G3").Value = "2016-08-28 05:05:41"

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2c621dc539b71fca51d393dd4715c8aea4fb2922/e2e/" + $newFile1, [Type]::Missing, [Type]::Missing, "e2e\" + $newFile1) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2c621dc539b71fca51d393dd4715c8aea4fb2922/e2e/" + $newFile2, [Type]::Missing, [Type]::Missing, "e2e\" + $newFile2) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newFile1
$wsZhCn.Range("G2").Value = "e1ef11d8-2e75-4ed4-b1d4-b8e55d6bb396.49819bff66962a9e3321e25b8aa4e76d9a5ccadb.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-28 05:05:36"
$wsZhCn.Range("I2").Value = $newFile1
$wsZhCn.Range("J2").Value = "e1ef11d8-2e75-4ed4-b1d4-b8e55d6bb396.49819bff66962a9e3321e25b8aa4e76d9a5ccadb.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-28 05:06:09"

$wsZhCn.Range("A3").Value = $newFile2
$wsZhCn.Range("G3").Value = "e1ef11d8-2e75-4ed4-b1d4-b8e55d6bb396.49819bff66962a9e3321e25b8aa4e76d9a5ccadb.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-28 05:05:36"
$wsZhCn.Range("I3").Value = $newFile2
$wsZhCn.Range("J3").Value = "e1ef11d8-2e75-4ed4-b1d4-b8e55d6bb396.49819bff66962a9e3321e25b8aa4e76d9a5ccadb.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-28 05:06:09"

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2c621dc539b71fca51d393dd4715c8aea4fb2922/e2e/" + $newFile1, [Type]::Missing, [Type]::Missing, $newFile1) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/6fb4f7426c8fa10f2b9139de5d401158144c0865/e2e/" + $newFile1, [Type]::Missing, [Type]::Missing, $newFile1) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2c621dc539b71fca51d393dd4715c8aea4fb2922/e2e/" + $newFile2, [Type]::Missing, [Type]::Missing, $newFile2) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/6fb4f7426c8fa10f2b9139de5d401158144c0865/e2e/" + $newFile2, [Type]::Missing, [Type]::Missing, $newFile2) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newFile1
$wsDeDe.Range("G2").Value = "e1ef11d8-2e75-4ed4-b1d4-b8e55d6bb396.49819bff66962a9e3321e25b8aa4e76d9a5ccadb.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-28 05:05:41"
$wsDeDe.Range("I2").Value = $newFile1
$wsDeDe.Range("J2").Value = "e1ef11d8-2e75-4ed4-b1d4-b8e55d6bb396.49819bff66962a9e3321e25b8aa4e76d9a5ccadb.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-28 05:06:16"

$wsDeDe.Range("A3").Value = $newFile2
$wsDeDe.Range("G3").Value = "e1ef11d8-2e75-4ed4-b1d4-b8e55d6bb396.49819bff66962a9e3321e25b8aa4e76d9a5ccadb.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-28 05:05:41"
$wsDeDe.Range("I3").Value = $newFile2
$wsDeDe.Range("J3").Value = "e1ef11d8-2e75-4ed4-b1d4-b8e55d6bb396.49819bff66962a9e3321e25b8aa4e76d9a5ccadb.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-28 05:06:16"

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2c621dc539b71fca51d393dd4715c8aea4fb2922/e2e/" + $newFile1, [Type]::Missing, [Type]::Missing, $newFile1) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d4d7b21dc016427dcd2fd0c9b147fd0ae7053df0/e2e/" + $newFile1, [Type]::Missing, [Type]::Missing, $newFile1) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2c621dc539b71fca51d393dd4715c8aea4fb2922/e2e/" + $newFile2, [Type]::Missing, [Type]::Missing, $newFile2) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d4d7b21dc016427dcd2fd0c9b147fd0ae7053df0/e2e/" + $newFile2, [Type]::Missing, [Type]::Missing, $newFile2) | Out-Null

"Handback status report regenerated."
